$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J - copy formatting from H1 (same header style)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J62
$data = @(
    @(7,9),
    @(7,7),
    @(6,7),
    @(7,8),
    @(6,7),
    @(7,7),
    @(7,8),
    @(8,8),
    @(7,7),
    @(7,8),
    @(5,6),
    @(5,5),
    @(7,7),
    @(8,8),
    @(10,11),
    @(8,8),
    @(7,7),
    @(9,9),
    @(7,7),
    @(8,8),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(8,8),
    @(6,7),
    @(6,7),
    @(8,9),
    @(8,8),
    @(8,8),
    @(7,8),
    @(6,7),
    @(8,9),
    @(8,8),
    @(7,8),
    @(7,8),
    @(7,8),
    @(7,8),
    @(7,8),
    @(4,5),
    @(6,6),
    @(6,7),
    @(8,8),
    @(8,8),
    @(6,6),
    @(7,7),
    @(4,6),
    @(10,10),
    @(5,7),
    @(7,8),
    @(8,9),
    @(9,9),
    @(8,9),
    @(7,8),
    @(7,7),
    @(8,8),
    @(5,5),
    @(8,8),
    @(4,4),
    @(2,2),
    @(6,6)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
